# settaggio timer e score funzionante
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update text "salvataggio gioco" -> "salvataggio gioco e load gioco"
# (it is used in cell D24)
$ws.Range("D24").Value = "salvataggio gioco e load gioco"

# Add "x" marks in F24 and G24
$ws.Range("F24").Value = "x"
$ws.Range("G24").Value = "x"

# Update the sheet view: scroll position and active selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("H26").Select()
